$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 17, shifting rows 17-40 down to 18-41
$ws.Rows.Item(17).Insert()

# Populate the new row 17 with data
$ws.Range("A17").Value = 3
$ws.Range("B17").Value = "Femacal de La Calera"
$ws.Range("C17").Value = "Coquimbo"
$ws.Range("D17").Value = 44526
$ws.Range("D17").NumberFormat = $ws.Range("D18").NumberFormat
$ws.Range("E17").Value = 5
$ws.Range("F17").Value = 100112022
$ws.Range("G17").Value = "Arveja Verde"
$ws.Range("H17").Value = "Perfection"
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 73
$ws.Range("K17").Value = 16000
$ws.Range("L17").Value = 17000
$ws.Range("M17").Value = 16521
$ws.Range("N17").Value = "`$/saco 25 kilos"
$ws.Range("O17").Value = "Provincia de Limarí"
$ws.Range("P17").Value = 661
$ws.Range("Q17").Value = 25
$ws.Range("R17").Value = "Hortaliza"
